$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.92773025893681
$ws.Range("D2").Value = 5.249664978213886
$ws.Range("E2").Value = 17.27897892724593
$ws.Range("F2").Value = 25.60861178551005
$ws.Range("G2").Value = 3.641342139503822
$ws.Range("I2").Value = 27.59913202289396
$ws.Range("K2").Value = 10.37788392761205
$ws.Range("L2").Value = 9.077777649411001
$ws.Range("M2").Value = 14.04963910846694
$ws.Range("O2").Value = 22.95529315108247
$ws.Range("B3").Value = 12.81108664042205
$ws.Range("D3").Value = 5.199758023703251
$ws.Range("E3").Value = 17.32460629892783
$ws.Range("F3").Value = 25.63610191954925
$ws.Range("G3").Value = 3.643167727173325
$ws.Range("I3").Value = 27.72514051326128
$ws.Range("K3").Value = 10.05590156360472
$ws.Range("L3").Value = 9.062470992886556
$ws.Range("M3").Value = 14.0226904257472
$ws.Range("O3").Value = 23.02378835839344
$ws.Range("B4").Value = 12.74094727003814
$ws.Range("D4").Value = 5.168426775786083
$ws.Range("E4").Value = 17.35421174792223
$ws.Range("F4").Value = 25.65991317009484
$ws.Range("G4").Value = 3.644348767417939
$ws.Range("I4").Value = 27.80698930560511
$ws.Range("K4").Value = 9.851071867117833
$ws.Range("L4").Value = 9.054393952642432
$ws.Range("M4").Value = 14.00793798179357
$ws.Range("O4").Value = 23.07097218498434
$ws.Range("B5").Value = 12.71276411819803
$ws.Range("D5").Value = 5.155491559604695
$ws.Range("E5").Value = 17.36667715919824
$ws.Range("F5").Value = 25.67135694637285
$ws.Range("G5").Value = 3.644845215809219
$ws.Range("I5").Value = 27.84147126864996
$ws.Range("K5").Value = 9.765886742457399
$ws.Range("L5").Value = 9.051437271060465
$ws.Range("M5").Value = 14.00238172905649
$ws.Range("O5").Value = 23.09148654904701
$ws.Range("B6").Value = 12.70810921358423
$ws.Range("D6").Value = 5.153333725810217
$ws.Range("E6").Value = 17.36877128057619
$ws.Range("F6").Value = 25.67336220390867
$ws.Range("G6").Value = 3.644928567984008
$ws.Range("I6").Value = 27.84726514380203
$ws.Range("K6").Value = 9.751640634771665
$ws.Range("L6").Value = 9.050966605705851
$ws.Range("M6").Value = 14.00148674891728
$ws.Range("O6").Value = 23.09497056505219
$ws.Range("B7").Value = 12.74056553066901
$ws.Range("D7").Value = 5.168252997644317
$ws.Range("E7").Value = 17.35437823588661
$ws.Range("F7").Value = 25.66006046182025
$ws.Range("G7").Value = 3.644355401230645
$ws.Range("I7").Value = 27.80744977236476
$ws.Range("K7").Value = 9.849929867857696
$ws.Range("L7").Value = 9.05435271912002
$ws.Range("M7").Value = 14.0078611984585
$ws.Range("O7").Value = 23.07124364339455
$ws.Range("B8").Value = 12.88722069255724
$ws.Range("D8").Value = 5.232603202743116
$ws.Range("E8").Value = 17.29438190152763
$ws.Range("F8").Value = 25.61665026141523
$ws.Range("G8").Value = 3.641959151974862
$ws.Range("I8").Value = 27.64165109165466
$ws.Range("K8").Value = 10.26838890641175
$ws.Range("L8").Value = 9.07222725127995
$ws.Range("M8").Value = 14.03997750522961
$ws.Range("O8").Value = 22.97784446372696
$ws.Range("B9").Value = 13.18534281472202
$ws.Range("D9").Value = 5.353094895045095
$ws.Range("E9").Value = 17.18929487216775
$ws.Range("F9").Value = 25.58660933386428
$ws.Range("G9").Value = 3.637735023729429
$ws.Range("I9").Value = 27.35198570439375
$ws.Range("K9").Value = 11.02935544807918
$ws.Range("L9").Value = 9.117644204676703
$ws.Range("M9").Value = 14.11699250719654
$ws.Range("O9").Value = 22.83549170763315
$ws.Range("B10").Value = 13.40913846286841
$ws.Range("D10").Value = 5.437835941959845
$ws.Range("E10").Value = 17.11967692103718
$ws.Range("F10").Value = 25.59817907132225
$ws.Range("G10").Value = 3.634918089829193
$ws.Range("I10").Value = 27.16068514359212
$ws.Range("K10").Value = 11.54849971937999
$ws.Range("L10").Value = 9.157159822242235
$ws.Range("M10").Value = 14.18184741273387
$ws.Range("O10").Value = 22.75593261603637
$ws.Range("B11").Value = 13.51162688828297
$ws.Range("D11").Value = 5.475501716839373
$ws.Range("E11").Value = 17.08963887982392
$ws.Range("F11").Value = 25.61074008544792
$ws.Range("G11").Value = 3.633698177179127
$ws.Range("I11").Value = 27.0783082261471
$ws.Range("K11").Value = 11.77537232153349
$ws.Range("L11").Value = 9.176430789693056
$ws.Range("M11").Value = 14.21308171450521
$ws.Range("O11").Value = 22.72520267467895
$ws.Range("B12").Value = 13.55050478049687
$ws.Range("D12").Value = 5.489632529522091
$ws.Range("E12").Value = 17.07849772906849
$ws.Range("F12").Value = 25.61654371119005
$ws.Range("G12").Value = 3.633245027162496
$ws.Range("I12").Value = 27.04778095269616
$ws.Range("K12").Value = 11.85990276150286
$ws.Range("L12").Value = 9.183910654093836
$ws.Range("M12").Value = 14.22515226777825
$ws.Range("O12").Value = 22.71435351358928
$ws.Range("B13").Value = 13.54212919958673
$ws.Range("D13").Value = 5.486595178625322
$ws.Range("E13").Value = 17.08088680156765
$ws.Range("F13").Value = 25.6152472779025
$ws.Range("G13").Value = 3.633342230223879
$ws.Range("I13").Value = 27.05432589590653
$ws.Range("K13").Value = 11.84175972062985
$ws.Range("L13").Value = 9.182291684273828
$ws.Range("M13").Value = 14.22254196062129
$ws.Range("O13").Value = 22.71665501333624
$ws.Range("B14").Value = 13.51482421152783
$ws.Range("D14").Value = 5.476666951023903
$ws.Range("E14").Value = 17.08871761461118
$ws.Range("F14").Value = 25.611196586235
$ws.Range("G14").Value = 3.633660720044298
$ws.Range("I14").Value = 27.0757833659332
$ws.Range("K14").Value = 11.78235461673136
$ws.Range("L14").Value = 9.177042530570326
$ws.Range("M14").Value = 14.21406993639548
$ws.Range("O14").Value = 22.72429431123827
$ws.Range("B15").Value = 13.49810704177878
$ws.Range("D15").Value = 5.470568225761795
$ws.Range("E15").Value = 17.09354460609251
$ws.Range("F15").Value = 25.60885168823098
$ws.Range("G15").Value = 3.633856949589457
$ws.Range("I15").Value = 27.08901352702708
$ws.Range("K15").Value = 11.745786130086
$ws.Range("L15").Value = 9.173850903854998
$ws.Range("M15").Value = 14.20891200902128
$ws.Range("O15").Value = 22.72907623177229
$ws.Range("B16").Value = 13.40245173581673
$ws.Range("D16").Value = 5.435356156950812
$ws.Range("E16").Value = 17.12167272167419
$ws.Range("F16").Value = 25.59750485107471
$ws.Range("G16").Value = 3.634999048362386
$ws.Range("I16").Value = 27.16616207589095
$ws.Range("K16").Value = 11.53348212687294
$ws.Range("L16").Value = 9.155926153947762
$ws.Range("M16").Value = 14.17984049660761
$ws.Range("O16").Value = 22.75805098922921
$ws.Range("B17").Value = 13.3439234549217
$ws.Range("D17").Value = 5.413524575679783
$ws.Range("E17").Value = 17.13934556037226
$ws.Range("F17").Value = 25.59241172453203
$ws.Range("G17").Value = 3.635715416390205
$ws.Range("I17").Value = 27.21467953101205
$ws.Range("K17").Value = 11.40082724782485
$ws.Range("L17").Value = 9.145258919842055
$ws.Range("M17").Value = 14.16244521960202
$ws.Range("O17").Value = 22.77722661195808
$ws.Range("B18").Value = 13.31032560405291
$ws.Range("D18").Value = 5.400884794602943
$ws.Range("E18").Value = 17.14966413515188
$ws.Range("F18").Value = 25.5901694114014
$ws.Range("G18").Value = 3.636133245862584
$ws.Range("I18").Value = 27.2430229173867
$ws.Range("K18").Value = 11.32365573232005
$ws.Range("L18").Value = 9.139245555322224
$ws.Range("M18").Value = 14.1526032777866
$ws.Range("O18").Value = 22.78876990465174
$ws.Range("B19").Value = 13.29896222602519
$ws.Range("D19").Value = 5.396591121418784
$ws.Range("E19").Value = 17.15318424339181
$ws.Range("F19").Value = 25.58952827191793
$ws.Range("G19").Value = 3.636275712094967
$ws.Range("I19").Value = 27.2526946645037
$ws.Range("K19").Value = 11.29737850590182
$ws.Range("L19").Value = 9.137230626869654
$ws.Range("M19").Value = 14.14929920715889
$ws.Range("O19").Value = 22.79276647179746
$ws.Range("B20").Value = 13.35014727345923
$ws.Range("D20").Value = 5.415857197063113
$ws.Range("E20").Value = 17.13744836420001
$ws.Range("F20").Value = 25.59288279636999
$ws.Range("G20").Value = 3.635638558492842
$ws.Range("I20").Value = 27.20946950095487
$ws.Range("K20").Value = 11.41503917588984
$ws.Range("L20").Value = 9.146381850838038
$ws.Range("M20").Value = 14.1642801132377
$ws.Range("O20").Value = 22.77513212243979
$ws.Range("B21").Value = 13.52284275103665
$ws.Range("D21").Value = 5.479586747486985
$ws.Range("E21").Value = 17.08641118303741
$ws.Range("F21").Value = 25.61235798021089
$ws.Range("G21").Value = 3.633566933252708
$ws.Range("I21").Value = 27.06946269114516
$ws.Range("K21").Value = 11.79984116008418
$ws.Range("L21").Value = 9.178579416644293
$ws.Range("M21").Value = 14.21655183623812
$ws.Range("O21").Value = 22.72202907267925
$ws.Range("B22").Value = 13.6360901972028
$ws.Range("D22").Value = 5.520463585337881
$ws.Range("E22").Value = 17.05441659494072
$ws.Range("F22").Value = 25.63118727391075
$ws.Range("G22").Value = 3.632264305761405
$ws.Range("I22").Value = 26.98184790595259
$ws.Range("K22").Value = 12.04326285783606
$ws.Range("L22").Value = 9.200683434667788
$ws.Range("M22").Value = 14.25212661819239
$ws.Range("O22").Value = 22.69191475506557
$ws.Range("B23").Value = 13.57562303211328
$ws.Range("D23").Value = 5.498719414994572
$ws.Range("E23").Value = 17.07136849242262
$ws.Range("F23").Value = 25.62058052328105
$ws.Range("G23").Value = 3.632954862706919
$ws.Range("I23").Value = 27.02825418873961
$ws.Range("K23").Value = 11.91409608159931
$ws.Range("L23").Value = 9.188790345686758
$ws.Range("M23").Value = 14.23301257862678
$ws.Range("O23").Value = 22.70756653731024
$ws.Range("B24").Value = 13.34733332443011
$ws.Range("D24").Value = 5.414802894225993
$ws.Range("E24").Value = 17.13830559367595
$ws.Range("F24").Value = 25.59266768842028
$ws.Range("G24").Value = 3.635673287317952
$ws.Range("I24").Value = 27.21182355355811
$ws.Range("K24").Value = 11.4086167831549
$ws.Range("L24").Value = 9.145873801617752
$ws.Range("M24").Value = 14.16345006260736
$ws.Range("O24").Value = 22.77607742474155
$ws.Range("B25").Value = 13.10373506951664
$ws.Range("D25").Value = 5.321141235872985
$ws.Range("E25").Value = 17.21638582211473
$ws.Range("F25").Value = 25.58882898896262
$ws.Range("G25").Value = 3.638827229840573
$ws.Range("I25").Value = 27.42656198591432
$ws.Range("K25").Value = 10.83027380762132
$ws.Range("L25").Value = 9.104264883609419
$ws.Range("M25").Value = 14.09468358437976
$ws.Range("O25").Value = 22.86961714590941
